# Adds three new AI/Nav-Mesh bullet steps and finishes the previously
# empty trailing list item under the "Artificial Intelligence (Nav Mesh)"
# heading, per the commit:
#   "added AI Nav Mesh, added MoveActor() in AI_Controller, and added a
#    RequestMove override function in movement component"

$d = $word.ActiveDocument

# Locate the heading paragraph, then the empty bulleted list paragraph
# that immediately follows it -- this is the paragraph that currently
# carries only the hidden "_GoBack" bookmark and no visible text.
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Artificial Intelligence (Nav Mesh)*") {
        $target = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the empty list paragraph after the Nav Mesh heading"
}

$r = $target.Range

# Preserve this paragraph's own identity (its w:rsidR/.../w:rsidP attributes)
# by reading them back off its OOXML rather than hard-coding them, so the
# first inserted sentence stays "inside" the original <w:p> element exactly
# like it would if someone had simply typed into it.
$oxml = $r.WordOpenXML
$origAttrs = ""
if ($oxml -match '<w:p\s+([^>]*?)/?>') {
    $origAttrs = [System.Text.RegularExpressions.Regex]::Replace($matches[1], 'w14:\w+="[^"]*"\s*', '').Trim()
    if ($origAttrs.Length -gt 0) {
        $origAttrs = " " + $origAttrs
    }
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr>'

# 1) "Place Nav Mesh Bounds Volume ..." -- stays inside the original <w:p>.
$para1 = "<w:p$origAttrs $wNs>$listPPr" +
    '<w:r><w:t xml:space="preserve">Place </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Nav</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Mesh Bounds Volume into the landscape (200x200) to start.  Make Z short.</w:t></w:r>' +
    '</w:p>'

# 2) "Create a Acceptance Radius Float member variable ..." -- new paragraph.
$para2 = "<w:p $wNs>$listPPr" +
    '<w:r><w:t xml:space="preserve">Create </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>a</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Acceptance Radius Float member variable in tank, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Uproperty</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'

# 3) "Call MoveActor on AI Controller ..." -- new paragraph.
$para3 = "<w:p $wNs>$listPPr" +
    '<w:r><w:t xml:space="preserve">Call </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>MoveActor</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> on AI Controller, pass in </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>AcceptanceRadius</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, and </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Playertank</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> in Tick</w:t></w:r>' +
    '</w:p>'

# 4) "Create a method in Movement Component, RequestDirectMove()" -- new
#    paragraph that now carries the bookmark (which moves to the new final
#    empty-looking-but-not-anymore paragraph, same as pressing Enter a few
#    times after typing text in Word would do).
$para4 = "<w:p $wNs>$listPPr" +
    '<w:r><w:t>Create a method in Movement Component, RequestDirectMove()</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$r.InsertXML($para1 + $para2 + $para3 + $para4) | Out-Null
Write-Host "Inserted 3 AI Nav Mesh steps and completed the final bullet."
